# Update "想去人数" (column F) values on the 展览 sheet and the 全部类型 sheet.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 298
$ws1.Range("F3").Value = 1197
$ws1.Range("F4").Value = 16817
$ws1.Range("F6").Value = 1646
$ws1.Range("F8").Value = 5
$ws1.Range("F10").Value = 219
$ws1.Range("F12").Value = 11660
$ws1.Range("F14").Value = 1339
$ws1.Range("F15").Value = 4625
$ws1.Range("F16").Value = 452

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 298
$ws4.Range("F4").Value = 1197
$ws4.Range("F5").Value = 16817
$ws4.Range("F7").Value = 1646
$ws4.Range("F9").Value = 5
$ws4.Range("F11").Value = 219
$ws4.Range("F15").Value = 11660
$ws4.Range("F17").Value = 1339
$ws4.Range("F18").Value = 4625
$ws4.Range("F19").Value = 452
